# Fill in the second (sample data) row of the resumenmesa spec sheet.
# Row 1 already holds the column headers; row 2 previously had empty
# placeholder cells (inline strings with no text). We now populate it
# with a concrete sample record.
#
# Numeric-looking values ("32", "1", "25", "180", "1207") must stay as
# *text* (matching the source dataset, which writes every column as a
# string), so they are entered with a leading apostrophe to force Excel
# to store them as text rather than auto-converting to numbers. The
# ClearFormats() calls afterwards strip the "number stored as text"
# quote-prefix flag that the apostrophe entry leaves on the cell style,
# since the target formatting is just the sheet's default style.
#
# The "resumen_votos_otro_estado" value (G2) is genuinely absent from
# the source record, so that cell is cleared instead of being given a
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "AR-X"
$ws.Range("C2").Value = "'32"
$ws.Range("D2").Value = "'1"
$ws.Range("E2").Value = "'25"
$ws.Range("F2").Value = "'180"
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = "'1207"

$ws.Range("C2:F2").ClearFormats()
$ws.Range("H2").ClearFormats()
